$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for A (date serial), B and C columns for rows 2-31 (shift date by 2 days,
# and replace B/C with newly retrained model output), and row 32 (date changes to a
# new time-of-day value, B/C stay 0).

$data = @(
    @{ Row = 2;  A = 46046;               B = 7.048;  C = 0.053 }
    @{ Row = 3;  A = 46046.01041666666;   B = 0.482;  C = 0.003 }
    @{ Row = 4;  A = 46046.02083333334;   B = 9.106;  C = 0.001 }
    @{ Row = 5;  A = 46046.03125;         B = 7.301;  C = 0.092 }
    @{ Row = 6;  A = 46046.04166666666;   B = 4.473;  C = 0 }
    @{ Row = 7;  A = 46046.05208333334;   B = 3.206;  C = 0 }
    @{ Row = 8;  A = 46046.0625;          B = 5.522;  C = 0.007 }
    @{ Row = 9;  A = 46046.07291666666;   B = 1.594;  C = 2.151 }
    @{ Row = 10; A = 46046.08333333334;   B = 1.749;  C = 0.348 }
    @{ Row = 11; A = 46046.09375;         B = 24.599; C = 0 }
    @{ Row = 12; A = 46046.10416666666;   B = 10.516; C = 0 }
    @{ Row = 13; A = 46046.11458333334;   B = 12.373; C = 0 }
    @{ Row = 14; A = 46046.125;           B = 4.932;  C = 0 }
    @{ Row = 15; A = 46046.13541666666;   B = 2.591;  C = 0 }
    @{ Row = 16; A = 46046.14583333334;   B = 2.34;   C = 0.015 }
    @{ Row = 17; A = 46046.15625;         B = 3.04;   C = 0.014 }
    @{ Row = 18; A = 46046.16666666666;   B = 3.222;  C = 0 }
    @{ Row = 19; A = 46046.17708333334;   B = 1.469;  C = 1.432 }
    @{ Row = 20; A = 46046.1875;          B = 2.578;  C = 0.945 }
    @{ Row = 21; A = 46046.19791666666;   B = 6.254;  C = 0.308 }
    @{ Row = 22; A = 46046.20833333334;   B = 2.324;  C = 0.985 }
    @{ Row = 23; A = 46046.21875;         B = 12.165; C = 0.475 }
    @{ Row = 24; A = 46046.22916666666;   B = 1.51;   C = 2.94 }
    @{ Row = 25; A = 46046.23958333334;   B = 6.293;  C = 0.067 }
    @{ Row = 26; A = 46046.25;            B = 1.112;  C = 0.09 }
    @{ Row = 27; A = 46046.26041666666;   B = 2.104;  C = 0.661 }
    @{ Row = 28; A = 46046.27083333334;   B = 9.829000000000001; C = 0.384 }
    @{ Row = 29; A = 46046.28125;         B = 11.851; C = 0.021 }
    @{ Row = 30; A = 46046.29166666666;   B = 16.62;  C = 0.263 }
    @{ Row = 31; A = 46046.30208333334;   B = 10.808; C = 0.042 }
    @{ Row = 32; A = 46046.3125;          B = 0;      C = 0 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
